# Path chose and letter attachment added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the displayed text of D5 (hyperlink cell) to a new email address,
# while keeping the existing hyperlink target (mailto:ivntz.apptest@mail.ru)
$ws.Range("D5").Value = "ivntz.apptest.main@mail.ru"

# Move/update the active selection to F13
$ws.Range("F13").Select()
